$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AIC_PinsFrames1")

# Test case row for "Invoice Lookup" - mark expected result as "pass"
$ws.Range("D2").Value = "pass"
